$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 37, pushing the existing rows 37-49 down to 39-51.
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# ---- New row 37 ----
$ws.Cells.Item(37,1).Value = 10
$ws.Cells.Item(37,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(37,3).Value = "La Araucanía"
$ws.Cells.Item(37,4).Value = 44559
$ws.Cells.Item(37,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37,5).Value = 9
$ws.Cells.Item(37,6).Value = 100112030
$ws.Cells.Item(37,7).Value = "Poroto granado"
$ws.Cells.Item(37,8).Value = "Sin especificar"
$ws.Cells.Item(37,9).Value = "Primera"
$ws.Cells.Item(37,10).Value = 25
$ws.Cells.Item(37,11).Value = 28000
$ws.Cells.Item(37,12).Value = 28000
$ws.Cells.Item(37,13).Value = 28000
$ws.Cells.Item(37,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(37,15).Value = "Región del Maule"
$ws.Cells.Item(37,16).Value = 1120
$ws.Cells.Item(37,17).Value = 25
$ws.Cells.Item(37,18).Value = "Hortaliza"

# ---- New row 38 ----
$ws.Cells.Item(38,1).Value = 10
$ws.Cells.Item(38,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38,3).Value = "La Araucanía"
$ws.Cells.Item(38,4).Value = 44559
$ws.Cells.Item(38,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38,5).Value = 9
$ws.Cells.Item(38,6).Value = 100112030
$ws.Cells.Item(38,7).Value = "Poroto granado"
$ws.Cells.Item(38,8).Value = "Sin especificar"
$ws.Cells.Item(38,9).Value = "Segunda"
$ws.Cells.Item(38,10).Value = 30
$ws.Cells.Item(38,11).Value = 25000
$ws.Cells.Item(38,12).Value = 25000
$ws.Cells.Item(38,13).Value = 25000
$ws.Cells.Item(38,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(38,15).Value = "Región del Maule"
$ws.Cells.Item(38,16).Value = 1000
$ws.Cells.Item(38,17).Value = 25
$ws.Cells.Item(38,18).Value = "Hortaliza"
